$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.596.28'
$ws.Range('E2').Value = '  +0.64%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.82'
$ws.Range('E3').Value = '  +0.78%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.88'
$ws.Range('E5').Value = '  -0.04%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.03%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4251'
$ws.Range('E7').Value = '  +0.74%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3643'
$ws.Range('E8').Value = '  +0.63%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.47'
$ws.Range('E9').Value = '  +0.27%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07299'
$ws.Range('E10').Value = '  +1.12%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8757'
$ws.Range('E11').Value = '  -3.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.55'
$ws.Range('E12').Value = '  -0.01%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.819.03'
$ws.Range('E13').Value = '  +2.94%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.329'
$ws.Range('E14').Value = '  -0.12%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.527'
$ws.Range('E15').Value = '  -0.60%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06879'
$ws.Range('E16').Value = '  +0.94%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  -0.01%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '79.83'
$ws.Range('E18').Value = '  +3.18%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008963'
$ws.Range('E19').Value = '  +0.24%  '

$ws.Range('E20').Value = '  +0.05%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.37'
$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.608.80'
$ws.Range('E22').Value = '  +0.67%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.991'
$ws.Range('E23').Value = '  +1.23%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.38'
$ws.Range('E24').Value = '  -1.52%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.053.58'
$ws.Range('E25').Value = '  +2.83%  '

$ws.Range('E26').Value = '  -1.92%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.42'
$ws.Range('E27').Value = '  +0.78%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.86'
$ws.Range('E28').Value = '  +3.88%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '122.24'
$ws.Range('E29').Value = '  +10.38%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.272'
$ws.Range('E30').Value = '  +0.14%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.880'
$ws.Range('E31').Value = '  +13.24%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08884'
$ws.Range('E32').Value = '  +0.42%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7676'
$ws.Range('E33').Value = '  -1.10%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.540'
$ws.Range('E34').Value = '  +0.69%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.972'
$ws.Range('E35').Value = '  +3.46%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.108'
$ws.Range('E36').Value = '  +3.69%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05380'
$ws.Range('E37').Value = '  +0.73%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.096'
$ws.Range('E38').Value = '  +2.13%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01939'
$ws.Range('E39').Value = '  +0.91%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.825'
$ws.Range('E40').Value = '  -3.92%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.905'
$ws.Range('E41').Value = '  +1.25%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5083'
$ws.Range('E42').Value = '  +0.61%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1652'
$ws.Range('E43').Value = '  +1.46%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.361'
$ws.Range('E44').Value = '  +1.68%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06529'
$ws.Range('E45').Value = '  -1.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.33'
$ws.Range('E46').Value = '  +0.99%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4689'
$ws.Range('E47').Value = '  -0.52%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '104.86'
$ws.Range('E48').Value = '  -0.10%  '

$ws.Range('E49').Value = '  +0.00%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.626'
$ws.Range('E50').Value = '  +0.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.41'
$ws.Range('E51').Value = '  +0.07%  '
